# Updated symbol list values per target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'311.06"
$ws.Range("E2").Value = "'0.92%"

$ws.Range("D3").Value = "'39.45"
$ws.Range("E3").Value = "'1.62%"

$ws.Range("D4").Value = "'5.122"
$ws.Range("E4").Value = "'0.32%"

$ws.Range("D5").Value = "'0.08108"
$ws.Range("E5").Value = "'-0.40%"

$ws.Range("D6").Value = "'2.018"
$ws.Range("E6").Value = "'2.88%"

$ws.Range("D7").Value = "'8.153"
$ws.Range("E7").Value = "'2.82%"

$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'4.237"
$ws.Range("E8").Value = "'0.82%"

$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9296"
$ws.Range("E9").Value = "'0.24%"

$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1429"
$ws.Range("E10").Value = "'-0.89%"

$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1932"
$ws.Range("E11").Value = "'-0.98%"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09086"
$ws.Range("E12").Value = "'-1.85%"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03500"
$ws.Range("E13").Value = "'-0.25%"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09818"
$ws.Range("E14").Value = "'-0.27%"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001402"
$ws.Range("E15").Value = "'0.00%"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005890"
$ws.Range("E16").Value = "'-2.38%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.791"
$ws.Range("E17").Value = "'5.32%"

$ws.Range("D18").Value = "'3.330"
$ws.Range("E18").Value = "'-4.18%"

$ws.Range("D19").Value = "'0.3453"
$ws.Range("E19").Value = "'0.20%"

$ws.Range("E20").Value = "'0.57%"

$ws.Range("D21").Value = "'4.727"
$ws.Range("E21").Value = "'-1.95%"

$ws.Range("D22").Value = "'0.2427"
$ws.Range("E22").Value = "'0.68%"

$ws.Range("D23").Value = "'0.04376"
$ws.Range("E23").Value = "'-1.60%"

$ws.Range("D24").Value = "'0.001230"
$ws.Range("E24").Value = "'-0.99%"

$ws.Range("D25").Value = "'0.004797"
$ws.Range("E25").Value = "'-1.18%"

$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'-0.05%"

$ws.Range("D27").Value = "'0.0004006"
$ws.Range("E27").Value = "'-9.92%"

$ws.Range("D39").Value = "'0.02128"
$ws.Range("E39").Value = "'1.27%"

$ws.Range("D40").Value = "'0.05101"
$ws.Range("E40").Value = "'-0.17%"

$ws.Range("D41").Value = "'0.007449"
$ws.Range("E41").Value = "'-0.31%"

$ws.Range("D42").Value = "'0.009826"
$ws.Range("E42").Value = "'-3.03%"

$ws.Range("E43").Value = "'-0.21%"

$ws.Range("E44").Value = "'-0.52%"

$ws.Range("D45").Value = "'0.008621"
$ws.Range("E45").Value = "'-17.95%"

$ws.Range("D46").Value = "'0.00006391"
$ws.Range("E46").Value = "'2.80%"

$ws.Range("E47").Value = "'-0.01%"

$ws.Range("E48").Value = "'-37.51%"

$ws.Range("D49").Value = "'0.002577"
$ws.Range("E49").Value = "'-16.07%"

$ws.Range("E50").Value = "'-0.01%"

$ws.Range("E51").Value = "'-0.01%"
